$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Dep or Non-Dep Without LOSOCV")
$sheet2 = $wb.Worksheets.Item("Dep or Non-Dep With LOSOCV")

# Sheet 1 - "Dep or Non-Dep Without LOSOCV"
# Row 10 (AdaBoost, first table)
$sheet1.Range("B10").Value = 0.87951807228915602
$sheet1.Range("C10").Value = 0.875
$sheet1.Range("D10").Value = 0.82352941176470495
$sheet1.Range("E10").Value = 0.84848484848484795
$sheet1.Range("F10").Value = 0.87094837935173997

# Row 20 (AdaBoost, second table)
$sheet1.Range("B20").Value = 0.69879518072289104
$sheet1.Range("C20").Value = 0.65517241379310298
$sheet1.Range("D20").Value = 0.55882352941176405
$sheet1.Range("E20").Value = 0.60317460317460303
$sheet1.Range("F20").Value = 0.67737094837935097
$sheet1.Range("F20").NumberFormat = "0.000000"

# Sheet 2 - "Dep or Non-Dep With LOSOCV"
# Row 10 (AdaBoost, first table)
$sheet2.Range("B10").Value = 0.83612402748766301
$sheet2.Range("C10").Value = 0.381818181818181
$sheet2.Range("D10").Value = 0.33167847304210901
$sheet2.Range("E10").Value = 0.34769674769674702
$sheet2.Range("F10").Value = 0.83612402748766301

# Row 20 (AdaBoost, second table)
$sheet2.Range("B20").Value = 0.70496024177842298
$sheet2.Range("C20").Value = 0.4
$sheet2.Range("D20").Value = 0.22222353404171499
$sheet2.Range("E20").Value = 0.26536998355180103
$sheet2.Range("F20").Value = 0.70496024177842298
